$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.0043348856584402963
$ws.Range("D2").Value = 0.19258042605313697
$ws.Range("E2").Value = 0.20957293835071567
$ws.Range("C3").Value = 0.0021986794029192937
$ws.Range("D3").Value = 0.066618869290621704
$ws.Range("E3").Value = 0.075237575684451483
$ws.Range("C4").Value = 0.0049363420186066039
$ws.Range("D4").Value = 0.3264275199460499
$ws.Range("E4").Value = 0.34577770794764251
$ws.Range("C5").Value = 0.0033877581710178373
$ws.Range("D5").Value = 0.16098830162245967
$ws.Range("E5").Value = 0.17426813358456053
$ws.Range("C6").Value = 0.0045618555684474229
$ws.Range("D6").Value = 0.360756187859059
$ws.Range("E6").Value = 0.37863840966476897
$ws.Range("C7").Value = 0.0037533146684885447
$ws.Range("D7").Value = 0.17729398808955057
$ws.Range("E7").Value = 0.19200678209144387
$ws.Range("C8").Value = 0.0044310348818369615
$ws.Range("D8").Value = 0.36529255039152031
$ws.Range("E8").Value = 0.3826619623329891
$ws.Range("C9").Value = 0.0036805931510827897
$ws.Range("D9").Value = 0.23238219037755944
$ws.Range("E9").Value = 0.24680991989656278
$ws.Range("C10").Value = 0.0043748068157648829
$ws.Range("D10").Value = 0.29341286317264087
$ws.Range("E10").Value = 0.31056186420146231
$ws.Range("C11").Value = 0.0039641853786518697
$ws.Range("D11").Value = 0.24463742595581925
$ws.Range("E11").Value = 0.26017682193321789
$ws.Range("C12").Value = 0.0031518324571125438
$ws.Range("D12").Value = 0.23527760910653236
$ws.Range("E12").Value = 0.24763261821343066
$ws.Range("C13").Value = 0.0044777253465341782
$ws.Range("D13").Value = 0.24962000866946965
$ws.Range("E13").Value = 0.26717245402496215
$ws.Range("C14").Value = 0.0040447064626449912
$ws.Range("D14").Value = 0.14024566775379421
$ws.Range("E14").Value = 0.15610069363498902
$ws.Range("C15").Value = 0.0046978060172384655
$ws.Range("D15").Value = 0.25098784595242679
$ws.Range("E15").Value = 0.26940299583921151
$ws.Range("C16").Value = 0.0035378334430600008
$ws.Range("D16").Value = 0.033417055825024258
$ws.Range("E16").Value = 0.047285167471967302
$ws.Range("C17").Value = 0.0056424034652907303
$ws.Range("D17").Value = 0.19645318463041225
$ws.Range("E17").Value = 0.21857110630571416
$ws.Range("C18").Value = 0.0043417201715328169
$ws.Range("D18").Value = -0.065327274136797495
$ws.Range("E18").Value = -0.048307970925473127
$ws.Range("C19").Value = 0.0052339738762650573
$ws.Range("D19").Value = 0.090026414480241845
$ws.Range("E19").Value = 0.11054331387567384
